$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.123541145015474, 2.651224325850143),
    @(3, 3.959010658874851, 2.48038747547068),
    @(4, 4.722695063536686, 3.580489688170352),
    @(5, 6.739021039846627, 2.457570659067509),
    @(6, 2.619839412265601, 1.929582042845213),
    @(7, -0.7919564768266385, 1.978648203842193),
    @(8, 1.877689851450803, 3.188468414048606),
    @(9, 2.705004599189187, 2.471779821159181),
    @(10, 1.110374544249249, 2.958768964947134),
    @(11, 2.267566233338814, 2.824054578526636),
    @(12, 2.688433258834588, 2.558605038029849),
    @(13, 1.014079695989589, 2.1961756525541),
    @(14, 3.013853578092252, 2.721358666668894),
    @(15, 1.331295149770684, 1.104792173470215),
    @(16, 0.04589006555719699, 1.367746832546346),
    @(17, 0.009546395482029624, 0.9083591402091473),
    @(18, 0.8709390141433015, 1.177606443599988),
    @(19, 0.7652063367885598, 1.680477107215861),
    @(20, 2.267579219134386, 2.225111592343887),
    @(21, 3.146753122914103, 1.799793194919874),
    @(22, 1.769033835366818, 0.2856860139923256),
    @(23, -4.774715709990263, 0.754926127539246),
    @(24, 1.95493704440024, 2.928378677701393),
    @(25, 3.478075069442799, 1.768431385360159),
    @(26, 1.232342134690434, 1.134779475590464),
    @(27, 0.2542811494408159, 1.353526127153426),
    @(28, 1.519778766382096, 0.7212678493511149),
    @(29, 1.469441753880329, 1.350400980499855),
    @(30, 1.638203081492495, 1.323745783269614),
    @(31, 2.268697431234346, 2.455413743911294),
    @(32, 1.984425467899631, 0.6687400825358569),
    @(33, 0.6066448776129052, 0.6230021429014077),
    @(34, -4.243076347305386, -2.163103471150829),
    @(35, 1.438499295329754, 1.882119284761474),
    @(36, 1.906593537051537, 1.473394465200051),
    @(37, 0.08348019664223827, 0.8702074629614476),
    @(38, -0.214505326882275, 0.8092352694139215),
    @(39, 0.1651547428133782, 0.7536567386490001)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
}
